$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the handoff/handback datetimes for the
# "4d8a515a..." row (row 3) to reflect the new report generation times.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-12 10:31:43"
$wsZh.Range("H3").Value = "2016-03-12 10:32:01"

# de-de sheet: same row, same kind of update.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-12 10:31:46"
$wsDe.Range("H3").Value = "2016-03-12 10:32:07"
